$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.002776666666666667
$ws.Range("H2").Value = 0.008330000000000001
$ws.Range("I2").Value = 0.2989949748743719
$ws.Range("J2").Value = 0.2989949748743719
$ws.Range("M2").Value = 0.002776666666666667
$ws.Range("N2").Value = 0.008330000000000001
$ws.Range("O2").Value = 0.2989949748743719
$ws.Range("P2").Value = 0.2989949748743719
$ws.Range("Q2").Value = [double]"7.709877777777779E-06"
$ws.Range("R2").Value = [double]"6.938890000000002E-05"
$ws.Range("S2").Value = 0.08939799500012631
$ws.Range("T2").Value = 0.08939799500012627
$ws.Range("G3").Value = 0.002776666666666667
$ws.Range("H3").Value = 0.008330000000000001
$ws.Range("I3").Value = 0.2989949748743719
$ws.Range("J3").Value = 0.2989949748743719
$ws.Range("O3").Value = 0.1584709260588658
$ws.Range("P3").Value = 0.1584709260588658
$ws.Range("Q3").Value = [double]"4.086327777777777E-06"
$ws.Range("R3").Value = [double]"3.677695E-05"
$ws.Range("S3").Value = 0.04738201055528902
$ws.Range("T3").Value = 0.04738201055528901
$ws.Range("G4").Value = 0.002776666666666667
$ws.Range("H4").Value = 0.008330000000000001
$ws.Range("I4").Value = 0.2989949748743719
$ws.Range("J4").Value = 0.2989949748743719
$ws.Range("M4").Value = 0.005038333333333333
$ws.Range("N4").Value = 0.015115
$ws.Range("O4").Value = 0.5425340990667624
$ws.Range("P4").Value = 0.5425340990667624
$ws.Range("Q4").Value = [double]"1.398977222222222E-05"
$ws.Range("R4").Value = 0.00012590795
$ws.Range("S4").Value = 0.1622149693189567
$ws.Range("T4").Value = 0.1622149693189566
$ws.Range("I5").Value = 0.1584709260588658
$ws.Range("J5").Value = 0.1584709260588658
$ws.Range("M5").Value = 0.002776666666666667
$ws.Range("N5").Value = 0.008330000000000001
$ws.Range("O5").Value = 0.2989949748743719
$ws.Range("P5").Value = 0.2989949748743719
$ws.Range("Q5").Value = [double]"4.086327777777777E-06"
$ws.Range("R5").Value = [double]"3.677695E-05"
$ws.Range("S5").Value = 0.04738201055528902
$ws.Range("T5").Value = 0.04738201055528901
$ws.Range("I6").Value = 0.1584709260588658
$ws.Range("J6").Value = 0.1584709260588658
$ws.Range("O6").Value = 0.1584709260588658
$ws.Range("P6").Value = 0.1584709260588658
$ws.Range("S6").Value = 0.0251130344059545
$ws.Range("T6").Value = 0.0251130344059545
$ws.Range("I7").Value = 0.1584709260588658
$ws.Range("J7").Value = 0.1584709260588658
$ws.Range("M7").Value = 0.005038333333333333
$ws.Range("N7").Value = 0.015115
$ws.Range("O7").Value = 0.5425340990667624
$ws.Range("P7").Value = 0.5425340990667624
$ws.Range("Q7").Value = [double]"7.414747222222221E-06"
$ws.Range("R7").Value = [double]"6.6732725E-05"
$ws.Range("S7").Value = 0.08597588109762226
$ws.Range("T7").Value = 0.08597588109762226
$ws.Range("G8").Value = 0.005038333333333333
$ws.Range("H8").Value = 0.015115
$ws.Range("I8").Value = 0.5425340990667624
$ws.Range("J8").Value = 0.5425340990667624
$ws.Range("M8").Value = 0.002776666666666667
$ws.Range("N8").Value = 0.008330000000000001
$ws.Range("O8").Value = 0.2989949748743719
$ws.Range("P8").Value = 0.2989949748743719
$ws.Range("Q8").Value = [double]"1.398977222222222E-05"
$ws.Range("R8").Value = 0.00012590795
$ws.Range("S8").Value = 0.1622149693189567
$ws.Range("T8").Value = 0.1622149693189566
$ws.Range("G9").Value = 0.005038333333333333
$ws.Range("H9").Value = 0.015115
$ws.Range("I9").Value = 0.5425340990667624
$ws.Range("J9").Value = 0.5425340990667624
$ws.Range("O9").Value = 0.1584709260588658
$ws.Range("P9").Value = 0.1584709260588658
$ws.Range("Q9").Value = [double]"7.414747222222221E-06"
$ws.Range("R9").Value = [double]"6.6732725E-05"
$ws.Range("S9").Value = 0.08597588109762226
$ws.Range("T9").Value = 0.08597588109762226
$ws.Range("G10").Value = 0.005038333333333333
$ws.Range("H10").Value = 0.015115
$ws.Range("I10").Value = 0.5425340990667624
$ws.Range("J10").Value = 0.5425340990667624
$ws.Range("M10").Value = 0.005038333333333333
$ws.Range("N10").Value = 0.015115
$ws.Range("O10").Value = 0.5425340990667624
$ws.Range("P10").Value = 0.5425340990667624
$ws.Range("Q10").Value = [double]"2.538480277777778E-05"
$ws.Range("R10").Value = 0.000228463225
$ws.Range("S10").Value = 0.2943432486501836
$ws.Range("T10").Value = 0.2943432486501836